# Update expected answer in test 5
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D2: "none" -> "x,left,right,temp,radius,radius1,quota,iter"
$ws.Range("D2").Value = "x,left,right,temp,radius,radius1,quota,iter"

# E2: (empty) -> "Get all variables on left hand side of assignment"
# Copy formatting from C2 (same look as other comment cells) then set the value.
$ws.Range("C2").Copy()
$ws.Range("E2").PasteSpecial(-4122)
$ws.Range("E2").Value = "Get all variables on left hand side of assignment"

# Update the view: scroll so column B is the left-most visible column,
# and move the active selection in the frozen (bottom-left) pane to E9.
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("E9").Select()
